$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; G=184.0626906666667; H=552.188072; I=0.6510505751503485; J=0.6510505751503486; K=3; M=29.57110033333333; N=88.713301; O=0.5311523066901919; P=0.5311523066901919; Q=5442.936293327297; R=48986.42663994568; S=0.3458070147630837; T=0.3458070147630838 }
    3 = @{ E=3; G=184.0626906666667; H=552.188072; I=0.6510505751503485; J=0.6510505751503486; K=3; M=12.180337; N=36.541011; O=0.2187816490048282; P=0.2187816490048282; Q=2241.945601446755; R=20177.51041302079; S=0.1424379184169351; T=0.1424379184169351 }
    4 = @{ E=3; G=184.0626906666667; H=552.188072; I=0.6510505751503485; J=0.6510505751503486; K=3; M=13.92204833333333; N=41.76614499999999; O=0.2500660443049799; P=0.2500660443049799; Q=2562.529675824715; R=23062.76708242244; S=0.1628056419703297; T=0.1628056419703297 }
    5 = @{ E=3; G=57.4434; H=172.3302; I=0.2031838091312023; J=0.2031838091312023; K=3; M=29.57110033333333; N=88.713301; O=0.5311523066901919; P=0.5311523066901919; Q=1698.6645448878; R=15287.9809039902; S=0.1079215489021378; T=0.1079215489021378 }
    6 = @{ E=3; G=57.4434; H=172.3302; I=0.2031838091312023; J=0.2031838091312023; K=3; M=12.180337; N=36.541011; O=0.2187816490048282; P=0.2187816490048282; Q=699.6799704258001; R=6297.119733832201; S=0.04445288881280671; T=0.04445288881280671 }
    7 = @{ E=3; G=57.4434; H=172.3302; I=0.2031838091312023; J=0.2031838091312023; K=3; M=13.92204833333333; N=41.76614499999999; O=0.2500660443049799; P=0.2500660443049799; Q=799.7297912309998; R=7197.568121078999; S=0.05080937141625782; T=0.05080937141625782 }
    8 = @{ E=3; G=41.21033366666666; H=123.631001; I=0.1457656157184491; J=0.1457656157184491; K=3; M=29.57110033333333; N=88.713301; O=0.5311523066901919; P=0.5311523066901919; Q=1218.634911627144; R=10967.7142046443; S=0.07742374302497033; T=0.07742374302497033 }
    9 = @{ E=3; G=41.21033366666666; H=123.631001; I=0.1457656157184491; J=0.1457656157184491; K=3; M=12.180337; N=36.541011; O=0.2187816490048282; P=0.2187816490048282; Q=501.9557519424457; R=4517.601767482011; S=0.03189084177508641; T=0.03189084177508641 }
    10 = @{ E=3; G=41.21033366666666; H=123.631001; I=0.1457656157184491; J=0.1457656157184491; K=3; M=13.92204833333333; N=41.76614499999999; O=0.2500660443049799; P=0.2500660443049799; Q=573.7322571401271; R=5163.590314261144; S=0.03645103091839238; T=0.03645103091839238 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
